# Add two new calculated columns to the cleaned payments worksheet:
#   G: TransactionSpeedNS = ResolveTime (D) - AttemptTime (C)
#   H: TransactionSpeedMS = ROUND(TransactionSpeedNS / 1,000,000)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("G1").Value = "TransactionSpeedNS"
$ws.Range("H1").Value = "TransactionSpeedMS"

# Find the last populated data row (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 51 }

for ($r = 2; $r -le $lastRow; $r++) {
    $attemptTime = $ws.Cells.Item($r, 3).Value2
    $resolveTime = $ws.Cells.Item($r, 4).Value2

    if ($attemptTime -ne $null -and $resolveTime -ne $null) {
        $speedNs = $resolveTime - $attemptTime
        $speedMs = [Math]::Round($speedNs / 1000000)

        $ws.Cells.Item($r, 7).Value = $speedNs
        $ws.Cells.Item($r, 8).Value = $speedMs
    }
}
